$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.153239130973816
$ws.Range("B1").Value = 2.373466014862061
$ws.Range("C1").Value = 5.083518505096436
$ws.Range("D1").Value = 2.283595561981201
$ws.Range("E1").Value = 1.243760108947754
